$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The fixture data set changes from a 4-column chart sample (Category/value1/value2/value3
# with Apple/Orange/Banana rows) to a 3-column upload/import test sample
# (Name/Id/TestFile with John Doe/Jane Doe rows).

# Drop the now-unused 4th column and 4th row first.
$ws.Columns.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Id"
$ws.Range("C1").Value = "TestFile"

# Data rows.
$ws.Range("A2").Value = "John Doe"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 123

$ws.Range("A3").Value = "Jane Doe"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 456

# Column widths (character units); Excel stores these with a small internal
# padding offset, so request slightly less than the target to land exactly
# on width 32 / 10 / 20 once persisted.
$ws.Columns.Item(1).ColumnWidth = 31.1666666667
$ws.Columns.Item(2).ColumnWidth = 9.1666666667
$ws.Columns.Item(3).ColumnWidth = 19.1666666667

# Keep the sheet tab color in the same red family as before.
$ws.Tab.Color = 252
